$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 177215.83
$ws.Range("I4").Value = 212459
$ws.Range("K4").Value = 212459
$ws.Range("M4").Value = -212345

$ws.Range("H17").Value = 1844
$ws.Range("J17").Value = 2726
$ws.Range("L17").Value = 8178
$ws.Range("N17").Value = -8514

$ws.Range("H33").Value = 1207.4706
$ws.Range("I33").Value = 219
$ws.Range("K33").Value = 219
$ws.Range("M33").Value = 10

$ws.Range("H107").Value = 6017.9473
$ws.Range("I107").Value = 5849.706
$ws.Range("K107").Value = 5849.706
$ws.Range("M107").Value = -3929.706

$ws.Range("H111").Value = 837.95
$ws.Range("I111").Value = 750.4375
$ws.Range("J111").Value = 1188
$ws.Range("K111").Value = 2251.3125
$ws.Range("L111").Value = 3564
$ws.Range("M111").Value = 815.6875
$ws.Range("N111").Value = -9698

$ws.Range("H121").Value = 2697.4167
$ws.Range("J121").Value = 2697.4167
$ws.Range("L121").Value = 8092.250100000001
$ws.Range("N121").Value = -11586.2501

$ws.Range("H138").Value = 3439.1836
$ws.Range("I138").Value = 1907.5
$ws.Range("J138").Value = 4051.8572
$ws.Range("K138").Value = 5722.5
$ws.Range("L138").Value = 12155.5716
$ws.Range("M138").Value = -582.5
$ws.Range("N138").Value = -22435.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1503.8572
$ws.Range("I2").Value = 1396.1724
$ws.Range("J2").Value = 2024.3334
$ws.Range("K2").Value = 1396.1724
$ws.Range("L2").Value = 2024.3334
$ws.Range("M2").Value = -1283.1724
$ws.Range("N2").Value = -2250.3334

$ws.Range("H110").Value = 3348.9333
$ws.Range("I110").Value = 3122.8147
$ws.Range("J110").Value = 5384
$ws.Range("K110").Value = 3122.8147
$ws.Range("L110").Value = 5384
$ws.Range("M110").Value = -1077.8147
$ws.Range("N110").Value = -9474

$ws.Range("H116").Value = 1503.8572
$ws.Range("I116").Value = 1396.1724
$ws.Range("J116").Value = 2024.3334
$ws.Range("K116").Value = 1396.1724
$ws.Range("L116").Value = 2024.3334
$ws.Range("M116").Value = 897.8276000000001
$ws.Range("N116").Value = -6612.3334

$ws.Range("H122").Value = 2251.7334
$ws.Range("I122").Value = 2168.7036
$ws.Range("K122").Value = 6506.110799999999
$ws.Range("M122").Value = -4056.110799999999

$ws.Range("H124").Value = 63427.062
$ws.Range("J124").Value = 63427.062
$ws.Range("L124").Value = 63427.062
$ws.Range("N124").Value = -73247.06200000001

$ws.Range("H138").Value = 153333.33
$ws.Range("J138").Value = 153333.33
$ws.Range("L138").Value = 153333.33
$ws.Range("N138").Value = -163613.33

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1503.8572
$ws.Range("I3").Value = 1396.1724
$ws.Range("J3").Value = 2024.3334
$ws.Range("K3").Value = 1396.1724
$ws.Range("L3").Value = 2024.3334
$ws.Range("M3").Value = -1282.1724
$ws.Range("N3").Value = -2252.3334

$ws.Range("H94").Value = 684.5417
$ws.Range("I94").Value = 554.0454999999999
$ws.Range("J94").Value = 2120
$ws.Range("K94").Value = 554.0454999999999
$ws.Range("L94").Value = 2120
$ws.Range("M94").Value = -103.0454999999999
$ws.Range("N94").Value = -3022

$ws.Range("H99").Value = 2722.3076
$ws.Range("I99").Value = 2775.5557
$ws.Range("K99").Value = 2775.5557
$ws.Range("M99").Value = -1277.5557

$ws.Range("H105").Value = 55558884
$ws.Range("I105").Value = 100002100
$ws.Range("K105").Value = 100002100
$ws.Range("M105").Value = -100000353

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1803.375
$ws.Range("J105").Value = 1650
$ws.Range("L105").Value = 1650
$ws.Range("N105").Value = -5144

$ws.Range("H134").Value = 2813.842
$ws.Range("I134").Value = 2351.3103
$ws.Range("J134").Value = 4304.222
$ws.Range("K134").Value = 7053.9309
$ws.Range("L134").Value = 12912.666
$ws.Range("M134").Value = -4518.9309
$ws.Range("N134").Value = -17982.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 2068.8
$ws.Range("I18").Value = 617.25
$ws.Range("K18").Value = 1851.75
$ws.Range("M18").Value = -1682.75

$ws.Range("H68").Value = 2066.682
$ws.Range("J68").Value = 2186.2307
$ws.Range("L68").Value = 6558.6921
$ws.Range("N68").Value = -8180.6921

$ws.Range("H71").Value = 2066.682
$ws.Range("J71").Value = 2186.2307
$ws.Range("L71").Value = 19676.0763
$ws.Range("N71").Value = -27788.0763

$ws.Range("H109").Value = 6404.8945
$ws.Range("I109").Value = 2268.4546
$ws.Range("K109").Value = 6805.3638
$ws.Range("M109").Value = -5765.3638

$ws.Range("H130").Value = 8753.25
$ws.Range("J130").Value = 10012.667
$ws.Range("L130").Value = 30038.001
$ws.Range("N130").Value = -40078.001

$ws.Range("H131").Value = 1797.8125
$ws.Range("J131").Value = 2058.16
$ws.Range("L131").Value = 6174.48
$ws.Range("N131").Value = -16254.48

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3202
$ws.Range("I80").Value = 1999
$ws.Range("J80").Value = 3502.75
$ws.Range("K80").Value = 1999
$ws.Range("L80").Value = 3502.75
$ws.Range("M80").Value = -1001
$ws.Range("N80").Value = -5498.75

$ws.Range("H83").Value = 3202
$ws.Range("I83").Value = 1999
$ws.Range("J83").Value = 3502.75
$ws.Range("K83").Value = 9995
$ws.Range("L83").Value = 17513.75
$ws.Range("M83").Value = -5003
$ws.Range("N83").Value = -27497.75

$ws.Range("H122").Value = 44733.69
$ws.Range("I122").Value = 61291.668
$ws.Range("J122").Value = 7478.25
$ws.Range("K122").Value = 183875.004
$ws.Range("L122").Value = 22434.75
$ws.Range("M122").Value = -181425.004
$ws.Range("N122").Value = -27334.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 19666
$ws.Range("J5").Value = 19666
$ws.Range("L5").Value = 19666
$ws.Range("N5").Value = -19892

$ws.Range("H46").Value = 1580.1875
$ws.Range("I46").Value = 721.4286
$ws.Range("J46").Value = 2248.111
$ws.Range("K46").Value = 721.4286
$ws.Range("L46").Value = 2248.111
$ws.Range("M46").Value = -533.4286
$ws.Range("N46").Value = -2624.111

$ws.Range("H55").Value = 1365.3462
$ws.Range("I55").Value = 284.6875
$ws.Range("K55").Value = 284.6875
$ws.Range("M55").Value = -111.6875

$ws.Range("H61").Value = 11313.412
$ws.Range("I61").Value = 9288.532999999999
$ws.Range("J61").Value = 26500
$ws.Range("K61").Value = 9288.532999999999
$ws.Range("L61").Value = 26500
$ws.Range("M61").Value = -9086.532999999999
$ws.Range("N61").Value = -26904

$ws.Range("H68").Value = 2391.2942
$ws.Range("I68").Value = 724.625
$ws.Range("J68").Value = 3872.7778
$ws.Range("K68").Value = 724.625
$ws.Range("L68").Value = 3872.7778
$ws.Range("M68").Value = 24.375
$ws.Range("N68").Value = -5370.7778

$ws.Range("H71").Value = 2391.2942
$ws.Range("I71").Value = 724.625
$ws.Range("J71").Value = 3872.7778
$ws.Range("K71").Value = 3623.125
$ws.Range("L71").Value = 19363.889
$ws.Range("M71").Value = 120.875
$ws.Range("N71").Value = -26851.889

$ws.Range("H93").Value = 1417.5555
$ws.Range("I93").Value = 1299
$ws.Range("J93").Value = 2899.5
$ws.Range("K93").Value = 1299
$ws.Range("L93").Value = 2899.5
$ws.Range("M93").Value = -51
$ws.Range("N93").Value = -5395.5

$ws.Range("H113").Value = 11313.412
$ws.Range("I113").Value = 9288.532999999999
$ws.Range("J113").Value = 26500
$ws.Range("K113").Value = 9288.532999999999
$ws.Range("L113").Value = 26500
$ws.Range("M113").Value = -7118.532999999999
$ws.Range("N113").Value = -30840

$ws.Range("H118").Value = 59000
$ws.Range("J118").Value = 59000
$ws.Range("L118").Value = 59000
$ws.Range("N118").Value = -62314

$ws.Range("H122").Value = 7546.343
$ws.Range("I122").Value = 7606.5293
$ws.Range("K122").Value = 22819.5879
$ws.Range("M122").Value = -20369.5879

$ws.Range("H136").Value = 4707.242
$ws.Range("I136").Value = 3743.238
$ws.Range("J136").Value = 6394.25
$ws.Range("K136").Value = 11229.714
$ws.Range("L136").Value = 19182.75
$ws.Range("M136").Value = -8679.714
$ws.Range("N136").Value = -24282.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H132").Value = 1207
$ws.Range("I132").Value = 999.4
$ws.Range("K132").Value = 2998.2
$ws.Range("M132").Value = -468.1999999999998
